$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-31 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-01 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("68÷2=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=5, 7", 2) | Out-Null
$d.Content.Find.Execute("64÷3=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "63÷8=7, 7", 2) | Out-Null
$d.Content.Find.Execute("49÷2=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=32, 1", 2) | Out-Null
$d.Content.Find.Execute("57÷7=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷7=13, 3", 2) | Out-Null
$d.Content.Find.Execute("35÷8=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2) | Out-Null
$d.Content.Find.Execute("91÷4=22, 3", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷9=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "45÷4=11, 1", 2) | Out-Null
$d.Content.Find.Execute("41÷9=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("97÷8=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=9, 1", 2) | Out-Null
$d.Content.Find.Execute("57÷4=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2) | Out-Null
$d.Content.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷6=1, 5", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=7, 0", 2) | Out-Null
$d.Content.Find.Execute("99÷3=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "31÷9=3, 4", 2) | Out-Null
$d.Content.Find.Execute("63÷9=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "23÷8=2, 7", 2) | Out-Null
$d.Content.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2) | Out-Null
$d.Content.Find.Execute("67÷4=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "38÷7=5, 3", 2) | Out-Null
$d.Content.Find.Execute("45÷6=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=10, 1", 2) | Out-Null
$d.Content.Find.Execute("38÷3=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "52÷3=17, 1", 2) | Out-Null
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=7, 4", 2) | Out-Null
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 2) | Out-Null
$d.Content.Find.Execute("60÷4=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=5, 6", 2) | Out-Null
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=11, 1", 2) | Out-Null
$d.Content.Find.Execute("70÷8=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=13, 3", 2) | Out-Null
$d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷2=39, 0", 2) | Out-Null
